$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (single-decimal-point numeric-looking strings) need to be forced to stay
# text, matching the source data which stores all Price/Volume cells as
# inline strings (not numbers).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "25.825.99"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "1.634.16"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  -1.73%  "

Set-TextValue "D5" "213.74"
$ws.Range("E5").Value = "  -1.23%  "

Set-TextValue "D6" "0.5006"
$ws.Range("E6").Value = "  -1.11%  "

Set-TextValue "D7" "1.002"
$ws.Range("E7").Value = "  -1.70%  "

Set-TextValue "D8" "0.2551"
$ws.Range("E8").Value = "  -1.27%  "

Set-TextValue "D9" "0.06319"
$ws.Range("E9").Value = "  -1.61%  "

Set-TextValue "D10" "19.25"
$ws.Range("E10").Value = "  -1.44%  "

Set-TextValue "D11" "0.07763"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "1.646.92"
$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.863.10"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "4.221"
$ws.Range("E14").Value = "  -0.91%  "

Set-TextValue "D15" "0.5383"
$ws.Range("E15").Value = "  -1.49%  "

$ws.Range("D16").Value = "0.0₅7811"
$ws.Range("E16").Value = "  -2.03%  "

Set-TextValue "D17" "64.06"
$ws.Range("E17").Value = "  +0.72%  "

$ws.Range("D18").Value = "25.878.49"
$ws.Range("E18").Value = "  -0.54%  "

Set-TextValue "D19" "1.002"
$ws.Range("E19").Value = "  -1.69%  "

Set-TextValue "D20" "194.87"
$ws.Range("E20").Value = "  -4.80%  "

Set-TextValue "D21" "4.331"
$ws.Range("E21").Value = "  +0.28%  "

Set-TextValue "D22" "9.804"
$ws.Range("E22").Value = "  -2.13%  "

Set-TextValue "D23" "5.928"
$ws.Range("E23").Value = "  -0.71%  "

Set-TextValue "D24" "1.004"
$ws.Range("E24").Value = "  -1.61%  "

Set-TextValue "D25" "1.890"
$ws.Range("E25").Value = "  -4.70%  "

Set-TextValue "D26" "139.55"
$ws.Range("E26").Value = "  -1.70%  "

Set-TextValue "D27" "0.1121"
$ws.Range("E27").Value = "  -3.07%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "15.57"
$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "6.753"
$ws.Range("E29").Value = "  -0.90%  "

Set-TextValue "D30" "1.234"
$ws.Range("E30").Value = "  -0.81%  "

Set-TextValue "D31" "0.04835"
$ws.Range("E31").Value = "  -3.14%  "

Set-TextValue "D32" "3.226"
$ws.Range("E32").Value = "  -1.19%  "

Set-TextValue "D33" "3.148"
$ws.Range("E33").Value = "  -1.90%  "

Set-TextValue "D34" "1.520"
$ws.Range("E34").Value = "  -1.20%  "

Set-TextValue "D35" "2.363"
$ws.Range("E35").Value = "  +0.25%  "

Set-TextValue "D36" "0.8813"
$ws.Range("E36").Value = "  -1.37%  "

Set-TextValue "D37" "2.588"
$ws.Range("E37").Value = "  -1.73%  "

$ws.Range("D38").Value = "1.122.45"
$ws.Range("E38").Value = "  +0.41%  "

Set-TextValue "D39" "0.5473"
$ws.Range("E39").Value = "  -3.36%  "

Set-TextValue "D40" "0.01554"
$ws.Range("E40").Value = "  -0.90%  "

Set-TextValue "D41" "1.001"
$ws.Range("E41").Value = "  -1.81%  "

Set-TextValue "D42" "5.652"
$ws.Range("E42").Value = "  +0.38%  "

Set-TextValue "D43" "0.8048"
$ws.Range("E43").Value = "  -1.52%  "

Set-TextValue "D44" "99.11"
$ws.Range("E44").Value = "  -0.64%  "

$ws.Range("D45").Value = "0.0₈120"
$ws.Range("E45").Value = "  +4.74%  "

$ws.Range("D46").Value = "1.775.00"
$ws.Range("E46").Value = "  -0.08%  "

Set-TextValue "D47" "0.4518"
$ws.Range("E47").Value = "  -1.18%  "

Set-TextValue "D48" "1.003"
$ws.Range("E48").Value = "  -1.22%  "

Set-TextValue "D49" "54.16"
$ws.Range("E49").Value = "  -1.22%  "

Set-TextValue "D50" "0.05059"

Set-TextValue "D51" "1.006"
$ws.Range("E51").Value = "  -1.23%  "
